$d = $word.ActiveDocument

# 1. Remove the "CompanyName" bookmark and the "Error? " run that precedes
#    "Communication System" in the title paragraph.
if ($d.Bookmarks.Exists("CompanyName")) {
    $d.Bookmarks("CompanyName").Delete()
}

$d.Content.Find.Execute("Error? Communication System", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Communication System", 2)

# 2. Remove the stray "test" run from the Revision History table cell.
$d.Content.Find.Execute("test", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
